$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (column D) and volume change (column E) values
# for rows 2-51, reflecting the latest scrape performed by the GitHub Actions
# workflow on Tue May  2 05:22:47 UTC 2023.

# The Price column stores plain-looking numbers (e.g. "1.002", "324.21") as
# text, not numeric values (note the "." thousands groupings used elsewhere
# in the same column, e.g. "27.979.28"). Force the column to Text first so
# Excel's automatic type detection does not silently convert these values
# into numbers when we assign them below.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.979.28"
$ws.Range("E2").Value = "  -2.38%  "
$ws.Range("D3").Value = "1.830.23"
$ws.Range("E3").Value = "  -1.21%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "324.21"
$ws.Range("E5").Value = "  -3.05%  "
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "0.4654"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").Value = "0.3867"
$ws.Range("E8").Value = "  -1.41%  "
$ws.Range("D9").Value = "0.07874"
$ws.Range("E9").Value = "  -0.29%  "
$ws.Range("D10").Value = "0.9593"
$ws.Range("E10").Value = "  -2.75%  "
$ws.Range("E11").Value = "  -1.83%  "
$ws.Range("D12").Value = "1.870.50"
$ws.Range("E12").Value = "  -6.27%  "
$ws.Range("D13").Value = "5.667"
$ws.Range("E13").Value = "  -3.35%  "
$ws.Range("D14").Value = "6.899"
$ws.Range("E14").Value = "  -1.92%  "
$ws.Range("D15").Value = "0.06810"
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("D16").Value = "87.23"
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").Value = "0.000009916"
$ws.Range("E18").Value = "  -1.62%  "
$ws.Range("D19").Value = "16.58"
$ws.Range("E19").Value = "  -2.68%  "
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").Value = "27.997.50"
$ws.Range("E21").Value = "  -2.39%  "
$ws.Range("D22").Value = "5.318"
$ws.Range("E22").Value = "  -1.51%  "
$ws.Range("D23").Value = "10.96"
$ws.Range("E23").Value = "  -2.40%  "
$ws.Range("D24").Value = "2.087"
$ws.Range("D25").Value = "2.052.27"
$ws.Range("E25").Value = "  -7.56%  "
$ws.Range("D26").Value = "153.64"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "19.13"
$ws.Range("E27").Value = "  -1.37%  "
$ws.Range("D28").Value = "5.723"
$ws.Range("E28").Value = "  -6.31%  "
$ws.Range("D29").Value = "1.965"
$ws.Range("E29").Value = "  -2.88%  "
$ws.Range("D30").Value = "117.39"
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("D31").Value = "0.09266"
$ws.Range("E31").Value = "  -1.66%  "
$ws.Range("D32").Value = "0.9314"
$ws.Range("E32").Value = "  -4.80%  "
$ws.Range("D33").Value = "5.281"
$ws.Range("E33").Value = "  -1.75%  "
$ws.Range("D34").Value = "1.316"
$ws.Range("E34").Value = "  -2.43%  "
$ws.Range("D35").Value = "3.290"
$ws.Range("E35").Value = "  -5.93%  "
$ws.Range("D36").Value = "0.05866"
$ws.Range("E36").Value = "  -4.32%  "
$ws.Range("D37").Value = "0.02140"
$ws.Range("E37").Value = "  -2.66%  "
$ws.Range("D38").Value = "1.144"
$ws.Range("E38").Value = "  -1.97%  "
$ws.Range("D39").Value = "7.770"
$ws.Range("E39").Value = "  +2.17%  "
$ws.Range("D40").Value = "0.5571"
$ws.Range("E40").Value = "  -2.42%  "
$ws.Range("D41").Value = "9.858"
$ws.Range("E41").Value = "  -2.40%  "
$ws.Range("D42").Value = "0.1758"
$ws.Range("E42").Value = "  -1.93%  "
$ws.Range("D43").Value = "11.61"
$ws.Range("E43").Value = "  -1.81%  "
$ws.Range("D44").Value = "0.5258"
$ws.Range("E44").Value = "  -2.35%  "
$ws.Range("D45").Value = "0.07008"
$ws.Range("E45").Value = "  -2.23%  "
$ws.Range("D46").Value = "2.144"
$ws.Range("E46").Value = "  -11.00%  "
$ws.Range("E47").Value = "  -4.44%  "
$ws.Range("D48").Value = "113.08"
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("E49").Value = "  -11.98%  "
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("D51").Value = "2.322"
$ws.Range("E51").Value = "  +0.91%  "
